$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.617.98'
$ws.Range("E2").Value = '  +2.29%  '

$ws.Range("D3").Value = '2.092.82'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E6").Value = '  +0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.25'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").Value = '2.397.24'
$ws.Range("E12").Value = '  +2.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.47%  '

$ws.Range("E15").Value = '  +1.55%  '

$ws.Range("E16").Value = '  +4.99%  '

$ws.Range("D17").Value = '2.093.32'
$ws.Range("E17").Value = '  +2.79%  '

$ws.Range("D18").Value = '38.554.09'
$ws.Range("E18").Value = '  +2.14%  '

$ws.Range("E19").Value = '  +3.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.23%  '

$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.81%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.54%  '

$ws.Range("E28").Value = '  +5.22%  '

$ws.Range("E29").Value = '  +1.98%  '

$ws.Range("E30").Value = '  +7.29%  '

$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("E32").Value = '  +5.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.00%  '

$ws.Range("E34").Value = '  +2.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0608'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("E36").Value = '  +1.72%  '

$ws.Range("E37").Value = '  +2.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.76%  '

$ws.Range("D41").Value = '1.546.28'
$ws.Range("E41").Value = '  +0.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0220'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.93%  '

$ws.Range("E44").Value = '  +1.01%  '

$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.74%  '

$ws.Range("E48").Value = '  +1.01%  '

$ws.Range("E49").Value = '  +2.66%  '

$ws.Range("E50").Value = '  +0.90%  '

$ws.Range("D51").Value = '2.288.41'
$ws.Range("E51").Value = '  +2.86%  '
